$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# commit message: "update notes through 2.1"
# The quiz note cells in column C previously embedded the due date in the
# text (e.g. "(2/4) Quiz 1 "); simplify them to just the quiz name since the
# dates live elsewhere / change.
$ws.Range("C16").Value = "Quiz 1 "
$ws.Range("C22").Value = "Quiz 2"
$ws.Range("C34").Value = "Quiz 3"
$ws.Range("C40").Value = "Quiz 4 "
$ws.Range("C46").Value = "Quiz 5"

# Scroll position / selection as last left by the author.
$ws.Range("E28").Select() | Out-Null
